$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 21: formatting-only change on G21 (centered, yellow fill) ---
$ws.Range("G21").HorizontalAlignment = -4108

# --- Row 27-38: add new "Comment" column C notes ---
$ws.Range("C27").Value = "Broderbund presents appears after 0:02"
$ws.Range("C28").Value = "A game by Jordan Mechner appears after 0:07"
$ws.Range("C29").Value = "Title screen appears after 0:15"
$ws.Range("C30").Value = "Sultan's absence screen"
$ws.Range("C31").Value = "Music continues from story_4_Jaffar_leaves"
$ws.Range("C33").Value = "Princess standing towards window"
$ws.Range("C34").Value = "Princess turns around"
$ws.Range("C35").Value = "Jaffar enters"
$ws.Range("C36").Value = "Hourglass appears after 0:13"
$ws.Range("C37").Value = "As sand starts to flow"
$ws.Range("C38").Value = "As Jaffar turns around and walks off"
$ws.Range("C27:C31").NumberFormat = "h:mm"
$ws.Range("C33:C38").NumberFormat = "h:mm"

# --- Row 30: G30 loses its highlight formatting (value unchanged) ---
$ws.Range("G30").Interior.ColorIndex = -4142

# --- Row 31: G31 "embrace" -> "n/a", add H31 filename ---
$ws.Range("G31").Value = "n/a"
$ws.Range("H31").Value = "pop_music_sumup"

# --- Row 33: add H33 filename ---
$ws.Range("H33").Value = "pop_music_princess"

# --- Row 37: G37 "n/a" -> "story_4_Jaffar_leaves", move audio/size/filename from row 38 to row 37 ---
$ws.Range("G37").Value = "story_4_Jaffar_leaves"
$ws.Range("H37").Value = "pop_music_leaves"
$ws.Range("I37").Value = 1
$ws.Range("J37").Value = "m-story4"

# --- Row 38: G38 "story_4_Jaffar_leaves" -> "n/a"; clear old I38/J38 ---
$ws.Range("G38").Value = "n/a"
$ws.Range("I38").ClearContents()
$ws.Range("J38").ClearContents()

# --- Sheet view state ---
$ws.Range("C37").Select()
$excel.ActiveWindow.ScrollRow = 4
